$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "last updated" timestamp in A1
$ws.Cells.Item(1,1).Value = "Datos actualizados a 31 de Marzo de 2020 a las 16:20"

# Update country names (column A) that shifted position due to re-sort
$ws.Cells.Item(73,1).Value = "Bulgaria"   # was "Letonia"
$ws.Cells.Item(74,1).Value = "Letonia"   # was "Bulgaria"
$ws.Cells.Item(105,1).Value = "Sri Lanka"   # was "Honduras"
$ws.Cells.Item(106,1).Value = "Honduras"   # was "Nigeria"
$ws.Cells.Item(107,1).Value = "Nigeria"   # was "Venezuela"
$ws.Cells.Item(108,1).Value = "Venezuela"   # was "Sri Lanka"
$ws.Cells.Item(132,1).Value = "Macao"   # was "Puerto Rico"
$ws.Cells.Item(133,1).Value = "Puerto Rico"   # was "Macao"
$ws.Cells.Item(158,1).Value = "Bahamas"   # was "Birmania"
$ws.Cells.Item(159,1).Value = "Birmania"   # was "Bahamas"
$ws.Cells.Item(167,1).Value = "Siria"   # was "Groenlandia"
$ws.Cells.Item(168,1).Value = "Groenlandia"   # was "Siria"
$ws.Cells.Item(170,1).Value = "Laos"   # was "Suazilandia"
$ws.Cells.Item(171,1).Value = "Suazilandia"   # was "Laos"
$ws.Cells.Item(173,1).Value = "Guinea-Bisau"   # was "Mozambique"
$ws.Cells.Item(174,1).Value = "Surinam"   # was "Libia"
$ws.Cells.Item(175,1).Value = "Mozambique"   # was "Guinea-Bisau"
$ws.Cells.Item(176,1).Value = "San Cristobal y Nieves"   # was "Surinam"
$ws.Cells.Item(177,1).Value = "Libia"   # was "San Cristobal y Nieves"
$ws.Cells.Item(182,1).Value = "Angola"   # was "Sudan"
$ws.Cells.Item(183,1).Value = "Sudan"   # was "Angola"
$ws.Cells.Item(184,1).Value = "Santa Sede"   # was "San Martin (Parte Holandesa)"
$ws.Cells.Item(185,1).Value = "San Martin (Parte Holandesa)"   # was "Santa Sede"
$ws.Cells.Item(186,1).Value = "Cabo Verde"   # was "Benin"
$ws.Cells.Item(187,1).Value = "Benin"   # was "San Bartolome"
$ws.Cells.Item(188,1).Value = "San Bartolome"   # was "Cabo Verde"
$ws.Cells.Item(190,1).Value = "Fiyi"   # was "Islas Turcas y Caicos"
$ws.Cells.Item(192,1).Value = "Islas Turcas y Caicos"   # was "Fiyi"
$ws.Cells.Item(197,1).Value = "Liberia"   # was "Belice"
$ws.Cells.Item(198,1).Value = "Belice"   # was "Republica de Africa Central"
$ws.Cells.Item(199,1).Value = "Islas Virgenes Britanicas"   # was "Liberia"
$ws.Cells.Item(200,1).Value = "Botsuana"   # was "Islas Virgenes Britanicas"
$ws.Cells.Item(201,1).Value = "Republica de Africa Central"   # was "Botsuana"
$ws.Cells.Item(205,1).Value = "Papua Nueva Guinea"   # was "Timor Oriental"
$ws.Cells.Item(206,1).Value = "Timor Oriental"   # was "Papua Nueva Guinea"

# Update numeric statistics (columns B-H) for rows whose data refreshed
$ws.Cells.Item(4,2).Value = 164800   # was 164665
$ws.Cells.Item(4,3).Value = 956   # was 821
$ws.Cells.Item(4,5).Value = 156115   # was 155981
$ws.Cells.Item(4,7).Value = 22   # was 21
$ws.Cells.Item(4,8).Value = 3178   # was 3177
$ws.Cells.Item(12,2).Value = 16186   # was 16176
$ws.Cells.Item(12,3).Value = 264   # was 254
$ws.Cells.Item(12,5).Value = 13968   # was 13980
$ws.Cells.Item(12,7).Value = 36   # was 14
$ws.Cells.Item(12,8).Value = 395   # was 373
$ws.Cells.Item(16,2).Value = 10019   # was 9974
$ws.Cells.Item(16,3).Value = 401   # was 356
$ws.Cells.Item(16,5).Value = 8796   # was 8751
$ws.Cells.Item(21,2).Value = 4681   # was 4667
$ws.Cells.Item(21,3).Value = 51   # was 37
$ws.Cells.Item(21,5).Value = 4387   # was 4374
$ws.Cells.Item(21,7).Value = 4   # was 3
$ws.Cells.Item(21,8).Value = 167   # was 166
$ws.Cells.Item(29,2).Value = 2738   # was 2449
$ws.Cells.Item(29,3).Value = 289   # was 0
$ws.Cells.Item(29,5).Value = 2570   # was 2285
$ws.Cells.Item(29,7).Value = 4   # was 0
$ws.Cells.Item(29,8).Value = 12   # was 8
$ws.Cells.Item(73,2).Value = 399   # was 398
$ws.Cells.Item(73,3).Value = 40   # was 22
$ws.Cells.Item(73,4).Value = 17   # was 1
$ws.Cells.Item(73,5).Value = 374   # was 397
$ws.Cells.Item(73,6).Value = 14   # was 3
$ws.Cells.Item(73,8).Value = 8   # was 0
$ws.Cells.Item(74,2).Value = 398   # was 379
$ws.Cells.Item(74,3).Value = 22   # was 20
$ws.Cells.Item(74,4).Value = 1   # was 17
$ws.Cells.Item(74,5).Value = 397   # was 354
$ws.Cells.Item(74,6).Value = 3   # was 13
$ws.Cells.Item(74,8).Value = 0   # was 8
$ws.Cells.Item(76,4).Value = 3   # was 7
$ws.Cells.Item(76,5).Value = 360   # was 356
$ws.Cells.Item(105,2).Value = 142   # was 141
$ws.Cells.Item(105,3).Value = 20   # was 2
$ws.Cells.Item(105,4).Value = 16   # was 3
$ws.Cells.Item(105,5).Value = 124   # was 131
$ws.Cells.Item(105,6).Value = 5   # was 4
$ws.Cells.Item(105,8).Value = 2   # was 7
$ws.Cells.Item(106,2).Value = 141   # was 135
$ws.Cells.Item(106,3).Value = 2   # was 4
$ws.Cells.Item(106,4).Value = 3   # was 8
$ws.Cells.Item(106,5).Value = 131   # was 125
$ws.Cells.Item(106,6).Value = 4   # was 0
$ws.Cells.Item(106,8).Value = 7   # was 2
$ws.Cells.Item(107,3).Value = 4   # was 0
$ws.Cells.Item(107,4).Value = 8   # was 39
$ws.Cells.Item(107,5).Value = 125   # was 93
$ws.Cells.Item(107,6).Value = 0   # was 6
$ws.Cells.Item(107,8).Value = 2   # was 3
$ws.Cells.Item(108,2).Value = 135   # was 132
$ws.Cells.Item(108,3).Value = 0   # was 10
$ws.Cells.Item(108,4).Value = 39   # was 16
$ws.Cells.Item(108,5).Value = 93   # was 114
$ws.Cells.Item(108,6).Value = 6   # was 5
$ws.Cells.Item(108,8).Value = 3   # was 2
$ws.Cells.Item(132,2).Value = 41   # was 39
$ws.Cells.Item(132,3).Value = 3   # was 0
$ws.Cells.Item(132,4).Value = 10   # was 1
$ws.Cells.Item(132,5).Value = 31   # was 36
$ws.Cells.Item(132,8).Value = 0   # was 2
$ws.Cells.Item(133,2).Value = 39   # was 38
$ws.Cells.Item(133,4).Value = 1   # was 10
$ws.Cells.Item(133,5).Value = 36   # was 28
$ws.Cells.Item(133,8).Value = 2   # was 0
$ws.Cells.Item(158,4).Value = 1   # was 0
$ws.Cells.Item(158,7).Value = 0   # was 1
$ws.Cells.Item(158,8).Value = 0   # was 1
$ws.Cells.Item(159,4).Value = 0   # was 1
$ws.Cells.Item(159,7).Value = 1   # was 0
$ws.Cells.Item(159,8).Value = 1   # was 0
$ws.Cells.Item(167,4).Value = 0   # was 2
$ws.Cells.Item(167,8).Value = 2   # was 0
$ws.Cells.Item(168,4).Value = 2   # was 0
$ws.Cells.Item(168,8).Value = 0   # was 2
$ws.Cells.Item(170,3).Value = 1   # was 0
$ws.Cells.Item(171,3).Value = 0   # was 1
$ws.Cells.Item(176,3).Value = 1   # was 0
$ws.Cells.Item(177,3).Value = 0   # was 1
$ws.Cells.Item(182,3).Value = 0   # was 1
$ws.Cells.Item(182,4).Value = 1   # was 0
$ws.Cells.Item(182,5).Value = 4   # was 5
$ws.Cells.Item(183,3).Value = 1   # was 0
$ws.Cells.Item(186,4).Value = 0   # was 1
$ws.Cells.Item(186,8).Value = 1   # was 0
$ws.Cells.Item(188,4).Value = 1   # was 0
$ws.Cells.Item(188,8).Value = 0   # was 1
$ws.Cells.Item(199,3).Value = 1   # was 0
$ws.Cells.Item(200,3).Value = 0   # was 1
